# Adds new HTTP-related error codes to the "Auth" sheet's error code table
# (rows 14 and 15, columns B/C) and makes "Auth" the active sheet/selection
# (previously "Register" was the active tab).

$wb = $excel.ActiveWorkbook

# --- Auth sheet: add the two new error rows ------------------------------
$authWs = $wb.Worksheets.Item("Auth")

$authWs.Range("B14").Value = "http return is null"
$authWs.Range("C14").Value = "ERROR"

$authWs.Range("B15").Value = "failed to parse http response"
$authWs.Range("C15").Value = "ERROR"

# --- Active sheet / selection changes -------------------------------------
# Previously "Register" was the selected/active tab with G10 selected; the
# active tab moves to "Auth" with B15 (the newly added row) selected.
$authWs.Select()
$authWs.Range("B15").Select()
